# Commit: "instead of showing figures through tk, generate pictures and export them"
#
# The spreadsheet is a static export (no live formulas - every cell is a
# cached value dropped in by the reporting script). On the
# "simple_exponential_smoothing" sheet the initial level used to be seeded
# with a pre-rounded figure (ROUND(AVERAGE(Demand),0) baked in at export
# time); the regenerated export now keeps full floating-point precision for
# that seed, which ripples through the whole Level/Forecast/Error/MAD/MAPE/TS
# column chain. Recompute that chain here with the same simple-exponential-
# smoothing model (alpha = 0.1) and drop the results back in as plain values
# -- matching how the rest of the workbook is already populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("simple_exponential_smoothing")

$alpha = 0.1
$lastRow = 14
$firstDataRow = 3

# Demand values live in column B, rows 3..14
$demand = @{}
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $demand[$r] = $ws.Cells.Item($r, 2).Value2
}

# Seed level (row 2, column C) = AVERAGE(Demand) at full precision (no rounding)
$sum = 0.0
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $sum += $demand[$r]
}
$n = $lastRow - $firstDataRow + 1
$level = @{}
$level[2] = $sum / $n
$ws.Cells.Item(2, 3).Value = $level[2]

$errSum = 0.0
$absErrSum = 0.0
$sqErrSum = 0.0
$pctErrAbsSum = 0.0
$count = 0

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $prevLevel = $level[$r - 1]

    $forecast = $prevLevel
    $levelNow = $alpha * $demand[$r] + (1 - $alpha) * $prevLevel
    $level[$r] = $levelNow

    $error = $forecast - $demand[$r]
    $absError = [Math]::Abs($error)

    $count += 1
    $errSum += $error
    $absErrSum += $absError
    $sqErrSum += [Math]::Pow($error, 2)
    $pctError = ($absError / $demand[$r]) * 100
    $pctErrAbsSum += $pctError

    $mse = $sqErrSum / $count
    $mad = $absErrSum / $count
    $mape = $pctErrAbsSum / $count
    $ts = $errSum / $mad

    $ws.Cells.Item($r, 3).Value = $levelNow     # C - Level
    $ws.Cells.Item($r, 4).Value = $forecast     # D - Forecast
    $ws.Cells.Item($r, 5).Value = $error        # E - Error
    $ws.Cells.Item($r, 6).Value = $absError     # F - Absolute Error
    $ws.Cells.Item($r, 7).Value = $mse          # G - Squared Error (MSE)
    $ws.Cells.Item($r, 8).Value = $mad          # H - MAD
    $ws.Cells.Item($r, 9).Value = $pctError     # I - % Error
    $ws.Cells.Item($r, 10).Value = $mape        # J - MAPE
    $ws.Cells.Item($r, 11).Value = $ts          # K - TS
}
